# The source workbook was re-opened and re-saved in a newer version of Excel.
# The meaningful, user-visible changes captured by that session were:
#   1. Column B was widened (to fit its header/content) and Column C was
#      given a custom width as well.
#   2. The user scrolled/clicked down into the data and left the selection
#      on cell C8.
#
# (Everything else in the diff - new xmlns/mc:Ignorable namespaces, the
# fileVersion/rupBuild bump, xr:revisionPtr / x15ac:absPath bookkeeping,
# the bookViews window geometry, the styles.xml x14ac/extLst slicer-style
# block, the theme's panose attributes and extra script font fallbacks,
# and the per-row x14ac:dyDescent hint - is metadata that Excel stamps on
# a file purely because it was opened/saved with a newer build; none of
# it is reachable through the Excel object model, and the cell data/
# formulas themselves are unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()

# --- Column widths -------------------------------------------------------
# Column B (the long descriptive-name column) and Column C get explicit
# custom widths.
$ws.Columns.Item(2).ColumnWidth = 40.5
$ws.Columns.Item(3).ColumnWidth = 9

# --- Selection -------------------------------------------------------
# Leave the final selection on C8, matching the saved sheet view.
$ws.Range("C8").Select() | Out-Null
